$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.84'
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.385'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06022'
$ws.Range("D5").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8148'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9280'
$ws.Range("D8").Style = "Normal"

$ws.Range("B9").Value = 'One'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01123'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '8OneONEBestin24h'

$ws.Range("B10").Value = 'WazirX'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1434'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '9WazirXWRX'

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07435'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03439'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03067'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09412'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '13BitMartTokenBMX'

$ws.Range("B15").Value = 'MCDex'

$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.007'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14MCDexMCB'

$ws.Range("B16").Value = 'BitForexToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001600'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '15BitForexTokenBF'

$ws.Range("B17").Value = 'CoinExToken'

$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04805'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '16CoinExTokenCET'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005651'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004158'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009877'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.667'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.423'
$ws.Range("D22").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1304'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.00007001'
$ws.Range("D26").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04021'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006411'
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002901'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.005953'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '43LocalTradersLCT'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005267'
$ws.Range("D45").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002515'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("D49").Style = "Normal"
